# Add a new review row (row 5) to the sheet, cloning the formatting of row 4
# (same app "com.singleton.stretchy" / "taxi game") and filling in the new
# reviewer's email addresses and review text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 4 into row 5 - this carries over both values and cell
# formatting (styles) in one shot.
$ws.Range("A4:F4").Copy($ws.Range("A5:F5"))
$ws.Rows.Item(5).RowHeight = $ws.Rows.Item(4).RowHeight

# Overwrite the cells that differ from row 4 with the new review's data.
$ws.Range("C5").Value2 = "shamirnaftali@gmail.com"
$ws.Range("D5").Value2 = "irisalmog47@gmail.com"
$ws.Range("F5").Value2 = "one of a kind taxi car game with great sound and graphics…"

# Match the author's final selection in the sheet.
$ws.Range("F5").Select()
